# "Fruta / hortaliza, semanal"
# Macroferia Regional de Talca - Pera: insert 4 new weekly price rows
# (Packham's Triumph, week of 2021-09-09) right before the existing
# tail rows, which shift down from 282-285 to 286-289.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 282..285; rows 282-285 (and everything after)
# move down to 286-289, same as EntireRow.Insert() repeated 4x.
$ws.Range("A282:A285").EntireRow.Insert()

$newRows = @(
    @{ Row=282; D=44448; K="Packham's Triumph"; L="Especial"; M=240; N=10000; O=10000; P=10000; Q="`$/bandeja 18 kilos granel"; R="Provincia de Linares"; S=556 },
    @{ Row=283; D=44448; K="Packham's Triumph"; L="Primera";  M=220; N=9000;  O=9000;  P=9000;  Q="`$/bandeja 18 kilos granel"; R="Provincia de Curicó";  S=500 },
    @{ Row=284; D=44448; K="Packham's Triumph"; L="Primera";  M=220; N=8000;  O=8000;  P=8000;  Q="`$/bandeja 18 kilos granel"; R="Provincia de Linares"; S=444 },
    @{ Row=285; D=44448; K="Packham's Triumph"; L="Segunda";  M=170; N=7000;  O=7000;  P=7000;  Q="`$/bandeja 18 kilos granel"; R="Provincia de Curicó";  S=389 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = 5
    $ws.Range("B$row").Value = "Macroferia Regional de Talca"
    $ws.Range("C$row").Value = "Maule"
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = 7
    $ws.Range("F$row").Value = "Fruta"
    $ws.Range("G$row").Value = 100104
    $ws.Range("H$row").Value = "Frutos de pepita"
    $ws.Range("I$row").Value = 100104005
    $ws.Range("J$row").Value = "Pera"
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = 18
}
